$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# The first three test cases (TestCase_E1..E3) had their Results recorded as
# PASS/FAIL; re-run them as SKIP (same value already used by the rest of the
# sheet), which also makes the now-unused "PASS"/"FAIL" shared strings drop
# out of the workbook on save.
$ws.Range("D2").Value = "SKIP"
$ws.Range("D3").Value = "SKIP"
$ws.Range("D4").Value = "SKIP"

# Leave the sheet's selection where the author left it after editing.
$ws.Range("C5").Select()
